$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Type" column for materials).
# Existing column B (E modulus / GPa) and everything to its right shifts
# right by one, carrying its formatting with it.
$ws.Columns.Item(2).Insert()

# Row 2 header text: "Material" -> "Label", and new "Type" header in B2.
$ws.Range("A2").Value = "Label"
$ws.Range("B2").Value = "Type"

# Existing material rows (3 and 4) get a "Steel" type in the new column B.
$ws.Range("B3").Value = "Steel"
$ws.Range("B4").Value = "Steel"

# New row 5: a timber section, styled like the plain "Normal" rows (same
# look as the sheet's default / non-bold rows).
$ws.Range("A5").Value = "100x100 C14"
$ws.Range("B5").Value = "Timber"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 350
$ws.Range("E5").Value = 14
$ws.Range("F5").Value = 8
$ws.Range("G5").Formula = "=10*10"
$ws.Range("H5").Formula = "=(10*10^3)/12"
$ws.Range("I5").Formula = "=H5"
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 5

# Normalise formatting so the (duplicate) bold style that used to live on
# the unit-row's middle columns collapses onto the same style as the rest
# of the bold header rows. Rows 3/4 already carry the correct plain style
# and are left untouched; only the brand-new row 5 needs it explicitly.
$ws.Range("A1:K2").Font.Bold = $true
$ws.Range("A5:K5").Font.Bold = $false

# Restore the selection to match the post-edit state.
$ws.Range("G16").Select() | Out-Null
